# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for the leves whose item prices moved since the last sync.
#
# Columns on every sheet:
#   H currentAveragePrice   I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ           L LevePriceHQ             M LeveProfitNQ
#   N LeveProfitHQ
#
# Only rows whose market data changed are touched; A:G (leve metadata) are
# left as-is.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H19").Value = 2310.75
$ws.Range("I19").Value = 932.3333
$ws.Range("J19").Value = 3137.8
$ws.Range("K19").Value = 932.3333
$ws.Range("L19").Value = 3137.8
$ws.Range("M19").Value = -757.3333
$ws.Range("N19").Value = -3487.8

$ws.Range("H40").Value = 2800
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325

$ws.Range("H41").Value = 313.2857
$ws.Range("I41").Value = 232.33333
$ws.Range("K41").Value = 232.33333
$ws.Range("M41").Value = 207.66667

$ws.Range("H43").Value = 1279.7059
$ws.Range("J43").Value = 1341.4286
$ws.Range("L43").Value = 1341.4286
$ws.Range("N43").Value = -1479.4286

$ws.Range("H55").Value = 306.36365
$ws.Range("J55").Value = 358
$ws.Range("L55").Value = 358
$ws.Range("N55").Value = -786

$ws.Range("H86").Value = 2996.3333
$ws.Range("I86").Value = 2999.5
$ws.Range("J86").Value = 2990
$ws.Range("K86").Value = 2999.5
$ws.Range("L86").Value = 2990
$ws.Range("M86").Value = -1876.5
$ws.Range("N86").Value = -5236

$ws.Range("H89").Value = 2996.3333
$ws.Range("I89").Value = 2999.5
$ws.Range("J89").Value = 2990
$ws.Range("K89").Value = 14997.5
$ws.Range("L89").Value = 14950
$ws.Range("M89").Value = -9381.5
$ws.Range("N89").Value = -26182

$ws.Range("H96").Value = 1563.75
$ws.Range("I96").Value = 1942
$ws.Range("K96").Value = 5826
$ws.Range("M96").Value = -4453

$ws.Range("H116").Value = 15173.375
$ws.Range("I116").Value = 34734.668
$ws.Range("J116").Value = 3436.6
$ws.Range("K116").Value = 34734.668
$ws.Range("L116").Value = 3436.6
$ws.Range("M116").Value = -31292.668
$ws.Range("N116").Value = -10320.6

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H124").Value = 46442
$ws.Range("J124").Value = 46442
$ws.Range("L124").Value = 46442
$ws.Range("N124").Value = -56262

$ws.Range("H132").Value = 925
$ws.Range("I132").Value = 932.36365
$ws.Range("J132").Value = 763
$ws.Range("K132").Value = 2797.09095
$ws.Range("L132").Value = 2289
$ws.Range("M132").Value = -267.0909499999998
$ws.Range("N132").Value = -7349

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 5702.1143
$ws.Range("I32").Value = 4081.9583
$ws.Range("K32").Value = 4081.9583
$ws.Range("M32").Value = -3794.9583

$ws.Range("H61").Value = 2266.4614
$ws.Range("I61").Value = 1373.174
$ws.Range("K61").Value = 1373.174
$ws.Range("M61").Value = -1161.174

$ws.Range("H74").Value = 388.14285
$ws.Range("I74").Value = 388.14285
$ws.Range("K74").Value = 388.14285
$ws.Range("M74").Value = 485.85715

$ws.Range("H77").Value = 388.14285
$ws.Range("I77").Value = 388.14285
$ws.Range("K77").Value = 1940.71425
$ws.Range("M77").Value = 2427.28575

$ws.Range("H136").Value = 2266.4614
$ws.Range("I136").Value = 1373.174
$ws.Range("K136").Value = 4119.522
$ws.Range("M136").Value = -1569.522

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H134").Value = 7634.1797
$ws.Range("I134").Value = 8763.286
$ws.Range("K134").Value = 26289.858
$ws.Range("M134").Value = -23754.858

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H60").Value = 11102.728
$ws.Range("J60").Value = 11102.728
$ws.Range("L60").Value = 11102.728
$ws.Range("N60").Value = -12124.728

$ws.Range("H99").Value = 1934.1666
$ws.Range("I99").Value = 1656.25
$ws.Range("J99").Value = 2490
$ws.Range("K99").Value = 1656.25
$ws.Range("L99").Value = 2490
$ws.Range("M99").Value = -158.25
$ws.Range("N99").Value = -5486

$ws.Range("H105").Value = 764.3333
$ws.Range("I105").Value = 877.2
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 877.2
$ws.Range("L105").Value = 200
$ws.Range("M105").Value = 869.8
$ws.Range("N105").Value = -3694

$ws.Range("H126").Value = 1934.1666
$ws.Range("I126").Value = 1656.25
$ws.Range("J126").Value = 2490
$ws.Range("K126").Value = 4968.75
$ws.Range("L126").Value = 7470
$ws.Range("M126").Value = -2498.75
$ws.Range("N126").Value = -12410

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H33").Value = 85.72727
$ws.Range("I33").Value = 112.28571
$ws.Range("K33").Value = 673.71426
$ws.Range("M33").Value = -390.71426

$ws.Range("H96").Value = 350
$ws.Range("I96").Value = 350
$ws.Range("K96").Value = 1050
$ws.Range("M96").Value = 1009

$ws.Range("H107").Value = 469.33334
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 496.92307
$ws.Range("K107").Value = 870
$ws.Range("L107").Value = 1490.76921
$ws.Range("M107").Value = 1050
$ws.Range("N107").Value = -5330.76921

$ws.Range("H108").Value = 2330.5
$ws.Range("I108").Value = 2330.5
$ws.Range("K108").Value = 6991.5
$ws.Range("M108").Value = -4111.5

$ws.Range("H122").Value = 883.4286
$ws.Range("J122").Value = 1097
$ws.Range("L122").Value = 9873
$ws.Range("N122").Value = -14773

$ws.Range("H130").Value = 125000810
$ws.Range("I130").Value = 125000810
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 375002430
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -374997410
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 10219102
$ws.Range("J131").Value = 15966.261
$ws.Range("L131").Value = 47898.783
$ws.Range("N131").Value = -57978.783

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H41").Value = 58358
$ws.Range("J41").Value = 58358
$ws.Range("L41").Value = 58358
$ws.Range("N41").Value = -59234

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 2698.75
$ws.Range("I132").Value = 2460.923
$ws.Range("J132").Value = 3729.3333
$ws.Range("K132").Value = 7382.768999999999
$ws.Range("L132").Value = 11187.9999
$ws.Range("M132").Value = -4852.768999999999
$ws.Range("N132").Value = -16247.9999
